$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("L3").Value = 1.37
$ws.Range("S3").Value = 3.1

# Row 4
$ws.Range("N4").Value = 1.34
$ws.Range("P4").Value = 1.34

# Row 8
$ws.Range("F8").Value = 1.52
$ws.Range("G8").Value = 1.54
$ws.Range("I8").Value = 7.2

# Row 9
$ws.Range("F9").Value = 3.75
$ws.Range("G9").Value = 4.3
$ws.Range("H9").Value = 1.89
$ws.Range("I9").Value = 2.08
$ws.Range("J9").Value = 3.65
$ws.Range("P9").Value = 2.04

# Row 10
$ws.Range("F10").Value = 2.06
$ws.Range("J10").Value = 3.45
$ws.Range("O10").Value = 1.48
$ws.Range("P10").Value = 1.66
$ws.Range("Q10").Value = 2.46
$ws.Range("R10").Value = 1.24
$ws.Range("S10").Value = 4.9
$ws.Range("T10").Value = 2.12
$ws.Range("U10").Value = 1.82
$ws.Range("Z10").Value = 30
$ws.Range("AJ10").Value = 26
$ws.Range("AN10").Value = 23
